$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force the COM host to keep a just-written run of text as its own
# <w:r> element (instead of silently re-merging it with a neighbouring run
# that happens to carry identical run properties) by toggling a character
# property on/off right after the text lands. The net formatting is
# unchanged, but the engine materialises a distinct run boundary.
# ---------------------------------------------------------------------------
function Split-Run($start, $end) {
    $r = $d.Range($start, $end)
    $r.Font.Bold = $true
    $r.Font.Bold = $false
}

# Appends " e único" (as two separate runs: " " and "e único") right after
# the word "Obrigatório" inside the given paragraph, turning
# "Obrigatório." into "Obrigatório e único."
function Add-EUnico($paraIndex) {
    $p = $d.Paragraphs($paraIndex)
    $scope = $p.Range.Duplicate
    $found = $scope.Find.Execute("Obrigatório", $false, $false, $false, $false, $false, `
        $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Add-EUnico: 'Obrigatório' not found in paragraph $paraIndex"
    }

    $scopeStart = $scope.Start
    $scopeEnd = $scope.End

    $pos = $scopeEnd
    $posAfter1 = $pos + 1
    $ins1 = $d.Range($pos, $pos)
    $ins1.InsertAfter(" ")
    Split-Run $pos $posAfter1

    $pos2 = $posAfter1
    $pos2After = $pos2 + 7
    $ins2 = $d.Range($pos2, $pos2)
    $ins2.InsertAfter("e único")
    Split-Run $pos2 $pos2After

    # Re-materialise the original "Obrigatório" run on its own so it is not
    # silently absorbed into the newly inserted text.
    Split-Run $scopeStart $scopeEnd
}

# ---------------------------------------------------------------------------
# 1) Title line "... GEOGRÁFICOS  -  IT575 ": the run holding "IT575" and the
#    run holding the trailing space collapse into the leading-space run, so
#    the paragraph keeps a single run reading " IT575 ".
# ---------------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(" IT575 ", $false, $false, $false, $false, $false, `
    $true, 1, $false, " IT575 ", 2)
if (-not $found1) {
    throw "' IT575 ' title run not found"
}

# ---------------------------------------------------------------------------
# 2) First "Obrigatório." (Objeto A / Propriedade A1, paragraph 28) becomes
#    "Obrigatório e único."
# ---------------------------------------------------------------------------
Add-EUnico 28

# ---------------------------------------------------------------------------
# 3) “teste” (paragraph 30) becomes “Olá Mundo!”, split across three runs:
#    “ / Olá Mundo! / ”
# ---------------------------------------------------------------------------
$p30 = $d.Paragraphs(30)
$rngT = $p30.Range.Duplicate
$foundT = $rngT.Find.Execute("“teste”", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
if (-not $foundT) {
    throw "'“teste”' not found in paragraph 30"
}
$tStart = $rngT.Start
$tEnd = $rngT.End

$newQuoted = [string]"“Olá Mundo!”"
$innerText = [string]"Olá Mundo!"
$innerLen = $innerText.Length

$replaceRng = $d.Range($tStart, $tEnd)
$replaceRng.Text = $newQuoted

$q1 = $tStart + 1
$q2 = $q1 + $innerLen
$q3 = $q2 + 1

Split-Run $tStart $q1
Split-Run $q1 $q2
Split-Run $q2 $q3

# ---------------------------------------------------------------------------
# 4) Second "Obrigatório." (Objeto B / Propriedade A1, paragraph 41) becomes
#    "Obrigatório e único."
# ---------------------------------------------------------------------------
Add-EUnico 41

Write-Output "done"
